# Plantilla Lista de Tareas de la Iteración
# Mockup CU 18, descripciones CU 17 y 19
#
# - Mockup para el CU 18 - Generar reporte de ingresos y egresos
# - Descripciones de los CU 17 - Registrar egreso y 19 - CRU Renta de espacio
# - Se actualizan las plantillas de casos de uso y lista de tareas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Row 24 - "Realizar mpckup de CU 20 y 21" (Mario): finished, 1 estimated hour,
# fully consumed on day 9 (column AI).
$ws.Range("G24").Value = 1
$ws.Range("AI24").Value = 1
$ws.Range("F24").Value = "Hecho"

# Row 26 - "Realizar descripciones de CU 14 y 16" (Mario): now marked done too.
$ws.Range("F26").Value = "Hecho"

# Row 27 - new task: mockup for CU 18 (Generar reporte), 2 estimated hours,
# fully consumed on day 10 (column AL); status done.
$ws.Range("D27").Value = "Realizar mockup de CU  18 Generar reporte"
$ws.Range("F27").Value = "Hecho"
$ws.Range("AL27").Value = 2

# Row 28 - new task: mockup for CU 22 (Iniciar sesión); status done.
$ws.Range("D28").Value = "Realizar mockup de CU 22 - Iniciar sesión"
$ws.Range("F28").Value = "Hecho"

# Row 29 - new task: descripción CU 17 y 19, 1 estimated hour, fully consumed
# on day 10 (column AL); status done.
$ws.Range("D29").Value = "Realizar descripción de CU 17 y 19"
$ws.Range("F29").Value = "Hecho"
$ws.Range("AL29").Value = 1

# Row 30 - new task: descripción CU 20 y 21; the hour previously logged on
# day 9 (column AI) is removed (task not finished yet) and its status moves
# back to "En proceso".
$ws.Range("D30").Value = "Realizar descripción de CU 20 y 21"
$ws.Range("AI30").ClearContents()
$ws.Range("F30").Value = "En proceso"

# Update the view: selection in the frozen bottom-right pane moves to H34.
$ws.Activate()
$ws.Range("H34").Select()
